$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I40").Value = 1000
$ws.Range("K40").Value = 1000
$ws.Range("M40").Value = -825

$ws.Range("H69").Value = 3113.5
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3113.5
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 9340.5
$ws.Range("N69").Value = -11088.5
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 3113.5
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3113.5
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 28021.5
$ws.Range("N72").Value = -36757.5
$ws.Range("M72").ClearContents()

$ws.Range("H116").Value = 2165
$ws.Range("I116").Value = 2118
$ws.Range("J116").Value = 2400
$ws.Range("K116").Value = 2118
$ws.Range("L116").Value = 2400
$ws.Range("M116").Value = 1324
$ws.Range("N116").Value = -9284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7291415.5
$ws.Range("I74").Value = 8359533.5
$ws.Range("J74").Value = 170627.33
$ws.Range("K74").Value = 8359533.5
$ws.Range("L74").Value = 170627.33
$ws.Range("M74").Value = -8358659.5
$ws.Range("N74").Value = -172375.33

$ws.Range("H77").Value = 7291415.5
$ws.Range("I77").Value = 8359533.5
$ws.Range("J77").Value = 170627.33
$ws.Range("K77").Value = 41797667.5
$ws.Range("L77").Value = 853136.6499999999
$ws.Range("M77").Value = -41793299.5
$ws.Range("N77").Value = -861872.6499999999

$ws.Range("H122").Value = 1473.7333
$ws.Range("I122").Value = 1127.875
$ws.Range("J122").Value = 1869
$ws.Range("K122").Value = 3383.625
$ws.Range("L122").Value = 5607
$ws.Range("M122").Value = -933.625
$ws.Range("N122").Value = -10507

$ws.Range("H132").Value = 45656.562
$ws.Range("I132").Value = 32150.787
$ws.Range("J132").Value = 75369.266
$ws.Range("K132").Value = 96452.361
$ws.Range("L132").Value = 226107.798
$ws.Range("M132").Value = -93922.361
$ws.Range("N132").Value = -231167.798

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 15672.529
$ws.Range("I86").Value = 17923.785
$ws.Range("J86").Value = 5166.6665
$ws.Range("K86").Value = 17923.785
$ws.Range("L86").Value = 5166.6665
$ws.Range("M86").Value = -16800.785
$ws.Range("N86").Value = -7412.6665

$ws.Range("H89").Value = 15672.529
$ws.Range("I89").Value = 17923.785
$ws.Range("J89").Value = 5166.6665
$ws.Range("K89").Value = 89618.925
$ws.Range("L89").Value = 25833.3325
$ws.Range("M89").Value = -84002.925
$ws.Range("N89").Value = -37065.3325

$ws.Range("H134").Value = 1822.8628
$ws.Range("I134").Value = 1757.4419
$ws.Range("J134").Value = 2174.5
$ws.Range("K134").Value = 5272.3257
$ws.Range("L134").Value = 6523.5
$ws.Range("M134").Value = -2737.3257
$ws.Range("N134").Value = -11593.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32306.768
$ws.Range("I31").Value = 24525.75
$ws.Range("J31").Value = 38142.53
$ws.Range("K31").Value = 24525.75
$ws.Range("L31").Value = 38142.53
$ws.Range("M31").Value = -24230.75
$ws.Range("N31").Value = -38732.53

$ws.Range("H34").Value = 32306.768
$ws.Range("I34").Value = 24525.75
$ws.Range("J34").Value = 38142.53
$ws.Range("K34").Value = 24525.75
$ws.Range("L34").Value = 38142.53
$ws.Range("M34").Value = -24323.75
$ws.Range("N34").Value = -38546.53

$ws.Range("H132").Value = 26014.977
$ws.Range("I132").Value = 1314.3448
$ws.Range("K132").Value = 3943.0344
$ws.Range("M132").Value = -1413.0344

$ws.Range("H134").Value = 31313.143
$ws.Range("I134").Value = 2924.818
$ws.Range("J134").Value = 79354.92
$ws.Range("K134").Value = 8774.454000000002
$ws.Range("L134").Value = 238064.76
$ws.Range("M134").Value = -6239.454000000002
$ws.Range("N134").Value = -243134.76

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 656.7875
$ws.Range("I68").Value = 688.9032
$ws.Range("J68").Value = 546.1667
$ws.Range("K68").Value = 2066.7096
$ws.Range("L68").Value = 1638.5001
$ws.Range("M68").Value = -1255.7096
$ws.Range("N68").Value = -3260.5001

$ws.Range("H71").Value = 656.7875
$ws.Range("I71").Value = 688.9032
$ws.Range("J71").Value = 546.1667
$ws.Range("K71").Value = 6200.1288
$ws.Range("L71").Value = 4915.5003
$ws.Range("M71").Value = -2144.1288
$ws.Range("N71").Value = -13027.5003

$ws.Range("H113").Value = 534.4103
$ws.Range("I113").Value = 507
$ws.Range("J113").Value = 543.86206
$ws.Range("K113").Value = 1521
$ws.Range("L113").Value = 1631.58618
$ws.Range("M113").Value = 649
$ws.Range("N113").Value = -5971.58618

$ws.Range("H131").Value = 990.8570999999999
$ws.Range("I131").Value = 400
$ws.Range("J131").Value = 1053.0526
$ws.Range("K131").Value = 1200
$ws.Range("L131").Value = 3159.1578
$ws.Range("M131").Value = 3840
$ws.Range("N131").Value = -13239.1578

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 929.0323
$ws.Range("I93").Value = 913.4
$ws.Range("K93").Value = 913.4
$ws.Range("M93").Value = 334.6

$ws.Range("H100").Value = 1567.9
$ws.Range("I100").Value = 1299
$ws.Range("J100").Value = 1971.25
$ws.Range("K100").Value = 1299
$ws.Range("L100").Value = 1971.25
$ws.Range("M100").Value = -758
$ws.Range("N100").Value = -3053.25

$ws.Range("H136").Value = 29869.178
$ws.Range("I136").Value = 19277.967
$ws.Range("J136").Value = 74503.57000000001
$ws.Range("K136").Value = 57833.901
$ws.Range("L136").Value = 223510.71
$ws.Range("M136").Value = -55283.901
$ws.Range("N136").Value = -228610.71

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1733.5385
$ws.Range("I81").Value = 737.5454999999999
$ws.Range("J81").Value = 2463.9333
$ws.Range("K81").Value = 1475.091
$ws.Range("L81").Value = 4927.8666
$ws.Range("M81").Value = -414.0909999999999
$ws.Range("N81").Value = -7049.8666

$ws.Range("H84").Value = 1733.5385
$ws.Range("I84").Value = 737.5454999999999
$ws.Range("J84").Value = 2463.9333
$ws.Range("K84").Value = 7375.455
$ws.Range("L84").Value = 24639.333
$ws.Range("M84").Value = -2071.455
$ws.Range("N84").Value = -35247.333

$ws.Range("H113").Value = 796.913
$ws.Range("I113").Value = 792.2273
$ws.Range("J113").Value = 900
$ws.Range("K113").Value = 2376.6819
$ws.Range("L113").Value = 2700
$ws.Range("M113").Value = -206.6819
$ws.Range("N113").Value = -7040

$ws.Range("H126").Value = 958.7143
$ws.Range("I126").Value = 837.9091
$ws.Range("J126").Value = 1401.6666
$ws.Range("K126").Value = 2513.7273
$ws.Range("L126").Value = 4204.9998
$ws.Range("M126").Value = -43.72730000000001
$ws.Range("N126").Value = -9144.9998

$ws.Range("H132").Value = 106350.266
$ws.Range("I132").Value = 91305.63
$ws.Range("J132").Value = 127036.625
$ws.Range("K132").Value = 273916.89
$ws.Range("L132").Value = 381109.875
$ws.Range("M132").Value = -271386.89
$ws.Range("N132").Value = -386169.875

$ws.Range("H136").Value = 41728.53
$ws.Range("I136").Value = 31508.363
$ws.Range("J136").Value = 60465.5
$ws.Range("K136").Value = 94525.08900000001
$ws.Range("L136").Value = 181396.5
$ws.Range("M136").Value = -91975.08900000001
$ws.Range("N136").Value = -186496.5
